$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Fix D153:D157 from text to numeric BSE codes (existing rows) ---
$ws.Cells.Item(153, 4).Value = 500387
$ws.Cells.Item(154, 4).Value = 532689
$ws.Cells.Item(155, 4).Value = 511196
$ws.Cells.Item(156, 4).Value = 500253
$ws.Cells.Item(157, 4).Value = 539876

# --- Append new rows 158-172 ---
$ws.Cells.Item(158, 1).Value = 1
$ws.Cells.Item(158, 2).Value = "ULTRACEMCO"
$ws.Cells.Item(158, 3).Value = "Ultratech Cement Limited"
$ws.Cells.Item(158, 4).Value = "'532538"
$ws.Cells.Item(158, 5).Value = -1.48
$ws.Cells.Item(158, 6).Value = 11658.7
$ws.Cells.Item(158, 7).Value = 556056
$ws.Cells.Item(158, 8).Value = "day"
$ws.Cells.Item(158, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(159, 1).Value = 2
$ws.Cells.Item(159, 2).Value = "ALKEM"
$ws.Cells.Item(159, 3).Value = "Alkem Laboratories Limited"
$ws.Cells.Item(159, 4).Value = "'539523"
$ws.Cells.Item(159, 5).Value = -1.67
$ws.Cells.Item(159, 6).Value = 5294
$ws.Cells.Item(159, 7).Value = 229725
$ws.Cells.Item(159, 8).Value = "day"
$ws.Cells.Item(159, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(160, 1).Value = 3
$ws.Cells.Item(160, 2).Value = "EICHERMOT"
$ws.Cells.Item(160, 3).Value = "Eicher Motors Limited"
$ws.Cells.Item(160, 4).Value = "'505200"
$ws.Cells.Item(160, 5).Value = 0.71
$ws.Cells.Item(160, 6).Value = 4916.1
$ws.Cells.Item(160, 7).Value = 664542
$ws.Cells.Item(160, 8).Value = "day"
$ws.Cells.Item(160, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(161, 1).Value = 4
$ws.Cells.Item(161, 2).Value = "ACC"
$ws.Cells.Item(161, 3).Value = "Acc Limited"
$ws.Cells.Item(161, 4).Value = "'500410"
$ws.Cells.Item(161, 5).Value = 0.75
$ws.Cells.Item(161, 6).Value = 2715.85
$ws.Cells.Item(161, 7).Value = 446135
$ws.Cells.Item(161, 8).Value = "day"
$ws.Cells.Item(161, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(162, 1).Value = 5
$ws.Cells.Item(162, 2).Value = "ASTRAL"
$ws.Cells.Item(162, 3).Value = "Astral Poly Technik Limited"
$ws.Cells.Item(162, 4).Value = "'532830"
$ws.Cells.Item(162, 5).Value = 0.19
$ws.Cells.Item(162, 6).Value = 2266.85
$ws.Cells.Item(162, 7).Value = 210636
$ws.Cells.Item(162, 8).Value = "day"
$ws.Cells.Item(162, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(163, 1).Value = 6
$ws.Cells.Item(163, 2).Value = "TATACOMM"
$ws.Cells.Item(163, 3).Value = "Tata Communications Limited"
$ws.Cells.Item(163, 4).Value = "'500483"
$ws.Cells.Item(163, 5).Value = -1.07
$ws.Cells.Item(163, 6).Value = 1851.05
$ws.Cells.Item(163, 7).Value = 223941
$ws.Cells.Item(163, 8).Value = "day"
$ws.Cells.Item(163, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(164, 1).Value = 7
$ws.Cells.Item(164, 2).Value = "MGL"
$ws.Cells.Item(164, 3).Value = "Mahanagar Gas Limited"
$ws.Cells.Item(164, 4).Value = "'539957"
$ws.Cells.Item(164, 5).Value = -0.3
$ws.Cells.Item(164, 6).Value = 1736.95
$ws.Cells.Item(164, 7).Value = 337792
$ws.Cells.Item(164, 8).Value = "day"
$ws.Cells.Item(164, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(165, 1).Value = 8
$ws.Cells.Item(165, 2).Value = "VOLTAS"
$ws.Cells.Item(165, 3).Value = "Voltas Limited"
$ws.Cells.Item(165, 4).Value = "'500575"
$ws.Cells.Item(165, 5).Value = -0.6
$ws.Cells.Item(165, 6).Value = 1521.7
$ws.Cells.Item(165, 7).Value = 690678
$ws.Cells.Item(165, 8).Value = "day"
$ws.Cells.Item(165, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(166, 1).Value = 9
$ws.Cells.Item(166, 2).Value = "SBICARD"
$ws.Cells.Item(166, 3).Value = "SBI Cards & Payment Services Ltd"
$ws.Cells.Item(166, 4).Value = "'543066"
$ws.Cells.Item(166, 5).Value = -1.04
$ws.Cells.Item(166, 6).Value = 730.9
$ws.Cells.Item(166, 7).Value = 1032255
$ws.Cells.Item(166, 8).Value = "day"
$ws.Cells.Item(166, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(167, 1).Value = 10
$ws.Cells.Item(167, 2).Value = "AMBUJACEM"
$ws.Cells.Item(167, 3).Value = "Ambuja Cements Limited"
$ws.Cells.Item(167, 4).Value = "'500425"
$ws.Cells.Item(167, 5).Value = 0.17
$ws.Cells.Item(167, 6).Value = 685.35
$ws.Cells.Item(167, 7).Value = 2484902
$ws.Cells.Item(167, 8).Value = "day"
$ws.Cells.Item(167, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(168, 1).Value = 11
$ws.Cells.Item(168, 2).Value = "MARICO"
$ws.Cells.Item(168, 3).Value = "Marico Limited"
$ws.Cells.Item(168, 4).Value = "'531642"
$ws.Cells.Item(168, 5).Value = 2.21
$ws.Cells.Item(168, 6).Value = 667.35
$ws.Cells.Item(168, 7).Value = 3414973
$ws.Cells.Item(168, 8).Value = "day"
$ws.Cells.Item(168, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(169, 1).Value = 12
$ws.Cells.Item(169, 2).Value = "GUJGASLTD"
$ws.Cells.Item(169, 3).Value = "Gujarat Gas Limited"
$ws.Cells.Item(169, 4).Value = "'539336"
$ws.Cells.Item(169, 5).Value = -0.73
$ws.Cells.Item(169, 6).Value = 637.8
$ws.Cells.Item(169, 7).Value = 1356263
$ws.Cells.Item(169, 8).Value = "day"
$ws.Cells.Item(169, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(170, 1).Value = 13
$ws.Cells.Item(170, 2).Value = "UPL"
$ws.Cells.Item(170, 3).Value = "Upl Limited"
$ws.Cells.Item(170, 4).Value = "'512070"
$ws.Cells.Item(170, 5).Value = -1.37
$ws.Cells.Item(170, 6).Value = 557.3
$ws.Cells.Item(170, 7).Value = 1674706
$ws.Cells.Item(170, 8).Value = "day"
$ws.Cells.Item(170, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(171, 1).Value = 14
$ws.Cells.Item(171, 2).Value = "LTF"
$ws.Cells.Item(171, 3).Value = "L&T Finance Ltd"
$ws.Cells.Item(171, 4).Value = "'533519"
$ws.Cells.Item(171, 5).Value = -0.8100000000000001
$ws.Cells.Item(171, 6).Value = 184.35
$ws.Cells.Item(171, 7).Value = 12011376
$ws.Cells.Item(171, 8).Value = "day"
$ws.Cells.Item(171, 9).Value = "16/07/2024 11:34:38"

$ws.Cells.Item(172, 1).Value = 15
$ws.Cells.Item(172, 2).Value = "IEX"
$ws.Cells.Item(172, 3).Value = "Indian Energy Exchange Ltd"
$ws.Cells.Item(172, 4).Value = "'540750"
$ws.Cells.Item(172, 5).Value = -0.49
$ws.Cells.Item(172, 6).Value = 177.34
$ws.Cells.Item(172, 7).Value = 11969274
$ws.Cells.Item(172, 8).Value = "day"
$ws.Cells.Item(172, 9).Value = "16/07/2024 11:34:38"

Write-Output "edit applied"